# Auto-generated edit script
# Applies numeric cell updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
# as described by the commit 'chore: update Sheets via scheduled runner'.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (55 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 450.83334
$ws.Range("I28").Value = 392.2381
$ws.Range("K28").Value = 392.2381
$ws.Range("M28").Value = 92.76190000000003
$ws.Range("H40").Value = 3079.8
$ws.Range("J40").Value = 5749.857
$ws.Range("L40").Value = 5749.857
$ws.Range("N40").Value = -6099.857
$ws.Range("H62").Value = 5514
$ws.Range("I62").Value = 4920
$ws.Range("K62").Value = 4920
$ws.Range("M62").Value = -4296
$ws.Range("H65").Value = 5514
$ws.Range("I65").Value = 4920
$ws.Range("K65").Value = 24600
$ws.Range("M65").Value = -21480
$ws.Range("H76").Value = 9000
$ws.Range("I76").Value = 9500
$ws.Range("K76").Value = 9500
$ws.Range("M76").Value = -9185
$ws.Range("H79").Value = 9000
$ws.Range("I79").Value = 9500
$ws.Range("K79").Value = 9500
$ws.Range("M79").Value = -8408
$ws.Range("H92").Value = 78832.55499999999
$ws.Range("J92").Value = 182070
$ws.Range("L92").Value = 182070
$ws.Range("N92").Value = -184566
$ws.Range("H96").Value = 15426.368
$ws.Range("I96").Value = 33371.625
$ws.Range("J96").Value = 2375.2727
$ws.Range("K96").Value = 100114.875
$ws.Range("L96").Value = 7125.8181
$ws.Range("M96").Value = -98741.875
$ws.Range("N96").Value = -9871.8181
$ws.Range("H98").Value = 59294.29
$ws.Range("I98").Value = 33486.535
$ws.Range("K98").Value = 33486.535
$ws.Range("M98").Value = -31988.535
$ws.Range("H106").Value = 5643.1055
$ws.Range("I106").Value = 5500.222
$ws.Range("K106").Value = 5500.222
$ws.Range("M106").Value = -4869.222
$ws.Range("H107").Value = 1818.0714
$ws.Range("I107").Value = 1742.875
$ws.Range("K107").Value = 1742.875
$ws.Range("M107").Value = 177.125
$ws.Range("H122").Value = 59294.29
$ws.Range("I122").Value = 33486.535
$ws.Range("K122").Value = 100459.605
$ws.Range("M122").Value = -98009.60500000001
$ws.Range("H137").Value = 2110.0454
$ws.Range("I137").Value = 1565.3889
$ws.Range("K137").Value = 4696.1667
$ws.Range("M137").Value = -2146.1667

# ---- Sheet: ARM (28 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1058.091
$ws.Range("I2").Value = 1197.421
$ws.Range("K2").Value = 1197.421
$ws.Range("M2").Value = -1084.421
$ws.Range("H32").Value = 17895.432
$ws.Range("I32").Value = 18768.578
$ws.Range("K32").Value = 18768.578
$ws.Range("M32").Value = -18481.578
$ws.Range("H40").Value = 35248
$ws.Range("J40").Value = 36331.332
$ws.Range("L40").Value = 36331.332
$ws.Range("N40").Value = -36683.332
$ws.Range("H46").Value = 7937.7
$ws.Range("J46").Value = 7707.1665
$ws.Range("L46").Value = 7707.1665
$ws.Range("N46").Value = -8345.166499999999
$ws.Range("H74").Value = 153229.83
$ws.Range("I74").Value = 102542.555
$ws.Range("K74").Value = 102542.555
$ws.Range("M74").Value = -101668.555
$ws.Range("H77").Value = 153229.83
$ws.Range("I77").Value = 102542.555
$ws.Range("K77").Value = 512712.775
$ws.Range("M77").Value = -508344.775
$ws.Range("H116").Value = 1058.091
$ws.Range("I116").Value = 1197.421
$ws.Range("K116").Value = 1197.421
$ws.Range("M116").Value = 1096.579

# ---- Sheet: BSM (29 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1058.091
$ws.Range("I3").Value = 1197.421
$ws.Range("K3").Value = 1197.421
$ws.Range("M3").Value = -1083.421
$ws.Range("H38").Value = 354177.34
$ws.Range("I38").Value = 2000000
$ws.Range("J38").Value = 25012.8
$ws.Range("K38").Value = 2000000
$ws.Range("L38").Value = 25012.8
$ws.Range("M38").Value = -1999584
$ws.Range("N38").Value = -25844.8
$ws.Range("H86").Value = 131138
$ws.Range("I86").Value = 3121
$ws.Range("J86").Value = 301827.34
$ws.Range("K86").Value = 3121
$ws.Range("L86").Value = 301827.34
$ws.Range("M86").Value = -1998
$ws.Range("N86").Value = -304073.34
$ws.Range("H89").Value = 131138
$ws.Range("I89").Value = 3121
$ws.Range("J89").Value = 301827.34
$ws.Range("K89").Value = 15605
$ws.Range("L89").Value = 1509136.7
$ws.Range("M89").Value = -9989
$ws.Range("N89").Value = -1520368.7
$ws.Range("H122").Value = 50390
$ws.Range("J122").Value = 50390
$ws.Range("L122").Value = 50390
$ws.Range("N122").Value = -60190

# ---- Sheet: CRP (8 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2522.7917
$ws.Range("J31").Value = 3622.6428
$ws.Range("L31").Value = 3622.6428
$ws.Range("N31").Value = -4212.6428
$ws.Range("H34").Value = 2522.7917
$ws.Range("J34").Value = 3622.6428
$ws.Range("L34").Value = 3622.6428
$ws.Range("N34").Value = -4026.6428

# ---- Sheet: CUL (19 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 116.22222
$ws.Range("J41").Value = 133.16667
$ws.Range("L41").Value = 399.50001
$ws.Range("N41").Value = -1075.50001
$ws.Range("H122").Value = 83379.09
$ws.Range("I122").Value = 1312.25
$ws.Range("J122").Value = 130274.43
$ws.Range("K122").Value = 11810.25
$ws.Range("L122").Value = 1172469.87
$ws.Range("M122").Value = -9360.25
$ws.Range("N122").Value = -1177369.87
$ws.Range("H139").Value = 3412.9546
$ws.Range("I139").Value = 3617.4
$ws.Range("K139").Value = 10852.2
$ws.Range("M139").Value = -5712.200000000001
$ws.Range("H140").Value = 2926.2222
$ws.Range("I140").Value = 1673.9
$ws.Range("K140").Value = 5021.700000000001
$ws.Range("M140").Value = 158.2999999999993

# ---- Sheet: GSM (11 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 10000
$ws.Range("I33").Value = 10000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 10000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -9748
$ws.Range("N33").ClearContents()
$ws.Range("H132").Value = 1254905.5
$ws.Range("I132").Value = 1543270.6
$ws.Range("K132").Value = 4629811.800000001
$ws.Range("M132").Value = -4627281.800000001

# ---- Sheet: LTW (16 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 65717.16
$ws.Range("J7").Value = 4669.3335
$ws.Range("L7").Value = 4669.3335
$ws.Range("N7").Value = -4893.3335
$ws.Range("H39").Value = 25000
$ws.Range("I39").Value = 25000
$ws.Range("K39").Value = 25000
$ws.Range("M39").Value = -24540
$ws.Range("H122").Value = 4209.3335
$ws.Range("I122").Value = 4209.3335
$ws.Range("K122").Value = 12628.0005
$ws.Range("M122").Value = -10178.0005
$ws.Range("H126").Value = 65717.16
$ws.Range("J126").Value = 4669.3335
$ws.Range("L126").Value = 14008.0005
$ws.Range("N126").Value = -18948.0005

# ---- Sheet: WVR (9 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 36666.555
$ws.Range("H107").Value = 639
$ws.Range("J107").Value = 1794
$ws.Range("L107").Value = 5382
$ws.Range("N107").Value = -9222
$ws.Range("H126").Value = 23848.334
$ws.Range("I126").Value = 26579.375
$ws.Range("K126").Value = 79738.125
$ws.Range("M126").Value = -77268.125

